$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $c = $ws.Range($cellAddr)
    $c.NumberFormat = "@"
    $c.Value2 = $val
    $c.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "37.510.31"
Set-TextValue "E2" "  -0.04%  "

# Row 3
Set-TextValue "D3" "2.082.46"
Set-TextValue "E3" "  +0.13%  "

# Row 4
Set-TextValue "E4" "  -0.01%  "

# Row 5
Set-TextValue "D5" "234.01"
Set-TextValue "E5" "  -0.45%  "

# Row 6
Set-TextValue "E6" "  +1.02%  "

# Row 7
$ws.Range("B7").Value2 = "Solana"
$ws.Range("C7").Value2 = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue "D7" "57.95"
Set-TextValue "E7" "  -0.72%  "

# Row 8
$ws.Range("B8").Value2 = "USDC"
$ws.Range("C8").Value2 = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue "D8" "1.00"
Set-TextValue "E8" "  +0.00%  "

# Row 9
Set-TextValue "D9" "0.389"
Set-TextValue "E9" "  +0.80%  "

# Row 10
Set-TextValue "D10" "0.0778"
Set-TextValue "E10" "  +1.70%  "

# Row 11
Set-TextValue "E11" "  +1.18%  "

# Row 12
Set-TextValue "D12" "2.388.18"
Set-TextValue "E12" "  +0.06%  "

# Row 13
Set-TextValue "D13" "14.44"
Set-TextValue "E13" "  -1.16%  "

# Row 14
Set-TextValue "D14" "21.04"
Set-TextValue "E14" "  -0.42%  "

# Row 15
Set-TextValue "D15" "0.782"
Set-TextValue "E15" "  +0.01%  "

# Row 16
Set-TextValue "E16" "  +0.07%  "

# Row 17
Set-TextValue "D17" "2.082.47"
Set-TextValue "E17" "  +0.15%  "

# Row 18
Set-TextValue "D18" "37.443.57"
Set-TextValue "E18" "  -0.66%  "

# Row 19
Set-TextValue "D19" "6.19"
Set-TextValue "E19" "  -1.20%  "

# Row 20
Set-TextValue "D20" "69.60"
Set-TextValue "E20" "  -1.32%  "

# Row 21
Set-TextValue "D21" "0.0₃0820"
Set-TextValue "E21" "  +0.28%  "

# Row 22
Set-TextValue "D22" "226.00"
Set-TextValue "E22" "  -0.49%  "

# Row 23
Set-TextValue "E23" "  +0.00%  "

# Row 24
Set-TextValue "D24" "2.46"
Set-TextValue "E24" "  +3.19%  "

# Row 25
Set-TextValue "E25" "  -3.04%  "

# Row 26
Set-TextValue "D26" "168.35"
Set-TextValue "E26" "  +0.99%  "

# Row 27
Set-TextValue "D27" "8.89"
Set-TextValue "E27" "  -1.81%  "

# Row 28
Set-TextValue "E28" "  -4.93%  "

# Row 29
Set-TextValue "D29" "0.132"
Set-TextValue "E29" "  +3.70%  "

# Row 30
Set-TextValue "D30" "19.17"
Set-TextValue "E30" "  -0.74%  "

# Row 31
Set-TextValue "E31" "  -0.74%  "

# Row 32
Set-TextValue "D32" "4.62"
Set-TextValue "E32" "  +2.05%  "

# Row 33
Set-TextValue "E33" "  -1.30%  "

# Row 34
Set-TextValue "E34" "  +0.11%  "

# Row 35
Set-TextValue "D35" "2.54"
Set-TextValue "E35" "  -1.69%  "

# Row 36
Set-TextValue "D36" "3.50"
Set-TextValue "E36" "  +3.99%  "

# Row 37
Set-TextValue "D37" "1.79"
Set-TextValue "E37" "  +0.78%  "

# Row 38
Set-TextValue "E38" "  -0.09%  "

# Row 39
Set-TextValue "E39" "  -5.53%  "

# Row 40
Set-TextValue "D40" "2.94"
Set-TextValue "E40" "  -0.43%  "

# Row 41
Set-TextValue "D41" "0.0957"
Set-TextValue "E41" "  +0.18%  "

# Row 42
Set-TextValue "D42" "1.481.96"
Set-TextValue "E42" "  -0.02%  "

# Row 43
Set-TextValue "D43" "96.99"
Set-TextValue "E43" "  +1.00%  "

# Row 45
Set-TextValue "E45" "  -2.23%  "

# Row 46
Set-TextValue "D46" "4.14"
Set-TextValue "E46" "  -12.04%  "

# Row 47
Set-TextValue "E47" "  +0.19%  "

# Row 48
Set-TextValue "D48" "15.51"
Set-TextValue "E48" "  -2.51%  "

# Row 49
Set-TextValue "D49" "7.28"
Set-TextValue "E49" "  -0.29%  "

# Row 50
Set-TextValue "D50" "3.00"
Set-TextValue "E50" "  +2.20%  "

# Row 51
Set-TextValue "D51" "2.274.03"
Set-TextValue "E51" "  +0.01%  "
